$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Requerimientos")

$ws.Range("B4").Value = "x"
$ws.Range("B5").Value = "x"

$ws.Activate()
$ws.Range("B6").Select()
